$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column price strings that look like plain numbers need the cell
# pre-formatted as Text, otherwise Excel auto-converts them to a Number
# (losing the original text formatting / introducing float artifacts).
$ws.Range("D2").Value = "42.710.87"
$ws.Range("E2").Value = "  -0.19%  "

$ws.Range("D3").Value = "2.546.90"
$ws.Range("E3").Value = "  +0.53%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.91"
$ws.Range("E5").Value = "  +4.87%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.25"
$ws.Range("E6").Value = "  -2.47%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.578"
$ws.Range("E7").Value = "  +0.36%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.536"
$ws.Range("E9").Value = "  -1.58%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.47"
$ws.Range("E10").Value = "  +0.37%  "

$ws.Range("E11").Value = "  -0.86%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.71"
$ws.Range("E12").Value = "  +0.09%  "

$ws.Range("E13").Value = "  -0.47%  "

$ws.Range("D14").Value = "2.937.25"
$ws.Range("E14").Value = "  +0.37%  "

$ws.Range("E15").Value = "  +6.41%  "

$ws.Range("D16").Value = "2.584.55"
$ws.Range("E16").Value = "  +4.13%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.870"
$ws.Range("E17").Value = "  +0.25%  "

$ws.Range("D18").Value = "42.754.86"
$ws.Range("E18").Value = "  -0.14%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.11"
$ws.Range("E19").Value = "  -0.66%  "

$ws.Range("E20").Value = "  +1.40%  "

$ws.Range("E21").Value = "  -1.66%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.14"
$ws.Range("E22").Value = "  -0.50%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "253.03"
$ws.Range("E23").Value = "  -0.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.97"
$ws.Range("E24").Value = "  +1.68%  "

$ws.Range("E25").Value = "  -2.60%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.30"
$ws.Range("E26").Value = "  -1.46%  "

$ws.Range("E27").Value = "  +0.06%  "

$ws.Range("E28").Value = "  +3.81%  "

$ws.Range("E29").Value = "  +4.91%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.30"
$ws.Range("E30").Value = "  +1.40%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.98"
$ws.Range("E31").Value = "  -3.81%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "156.23"
$ws.Range("E32").Value = "  -0.60%  "

$ws.Range("E33").Value = "  +1.48%  "

$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.37"
$ws.Range("E34").Value = "  +2.35%  "

$ws.Range("B35").Value = "Celestia"
$ws.Range("C35").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.29"
$ws.Range("E35").Value = "  -0.56%  "

$ws.Range("E36").Value = "  -0.17%  "

$ws.Range("E37").Value = "  +0.04%  "

$ws.Range("E38").Value = "  -2.30%  "

$ws.Range("E39").Value = "  +13.84%  "

$ws.Range("E40").Value = "  -0.02%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.70"
$ws.Range("E41").Value = "  -4.93%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.86"
$ws.Range("E42").Value = "  +0.41%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.37"
$ws.Range("E43").Value = "  -1.18%  "

$ws.Range("E44").Value = "  +0.30%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0303"
$ws.Range("E45").Value = "  -0.60%  "

$ws.Range("D46").Value = "2.035.50"
$ws.Range("E46").Value = "  -2.79%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "84.70"
$ws.Range("E47").Value = "  -1.98%  "

$ws.Range("E48").Value = "  -0.02%  "

$ws.Range("D49").Value = "2.790.72"
$ws.Range("E49").Value = "  +0.30%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.12"
$ws.Range("E50").Value = "  +0.96%  "

$ws.Range("E51").Value = "  -0.29%  "
